$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddProduct")

$ws.Range("L4").Value = "gio-qua-2.jpg"
$ws.Range("L6").Select()
